{"js": "const replacements = [\n  [\"27\u00d753=1431\", \"72\u00d739=2808\"],\n  [\"24\u00d748=1152\", \"82\u00d734=2788\"],\n  [\"88\u00d731=2728\", \"94\u00d781=7614\"],\n  [\"38\u00d713=494\", \"35\u00d797=3395\"],\n  [\"51\u00d717=867\", \"33\u00d766=2178\"],\n  [\"43\u00d799=4257\", \"19\u00d778=1482\"],\n  [\"62\u00d764=3968\", \"58\u00d746=2668\"],\n  [\"51\u00d720=1020\", \"93\u00d769=6417\"],\n  [\"11\u00d724=264\", \"23\u00d788=2024\"],\n  [\"64\u00d728=1792\", \"52\u00d719=988\"],\n  [\"73\u00d735=2555\", \"50\u00d744=2200\"],\n  [\"66\u00d747=3102\", \"97\u00d759=5723\"],\n  [\"78\u00d798=7644\", \"97\u00d781=7857\"],\n  [\"79\u00d765=5135\", \"27\u00d794=2538\"],\n  [\"20\u00d722=440\", \"14\u00d755=770\"],\n  [\"58\u00d724=1392\", \"42\u00d746=1932\"],\n  [\"92\u00d760=5520\", \"72\u00d720=1440\"],\n  [\"34\u00d731=1054\", \"80\u00d734=2720\"],\n  [\"65\u00d789=5785\", \"82\u00d732=2624\"],\n  [\"51\u00d773=3723\", \"67\u00d755=3685\"],\n  [\"36\u00d741=1476\", \"26\u00d741=1066\"],\n  [\"60\u00d737=2220\", \"62\u00d798=6076\"],\n  [\"13\u00d715=195\", \"42\u00d775=3150\"],\n  [\"78\u00d720=1560\", \"17\u00d740=680\"],\n  [\"29\u00d787=2523\", \"50\u00d780=4000\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @{old=\"27\u00d753=1431\"; new=\"72\u00d739=2808\"},\n  @{old=\"24\u00d748=1152\"; new=\"82\u00d734=2788\"},\n  @{old=\"88\u00d731=2728\"; new=\"94\u00d781=7614\"},\n  @{old=\"38\u00d713=494\"; new=\"35\u00d797=3395\"},\n  @{old=\"51\u00d717=867\"; new=\"33\u00d766=2178\"},\n  @{old=\"43\u00d799=4257\"; new=\"19\u00d778=1482\"},\n  @{old=\"62\u00d764=3968\"; new=\"58\u00d746=2668\"},\n  @{old=\"51\u00d720=1020\"; new=\"93\u00d769=6417\"},\n  @{old=\"11\u00d724=264\"; new=\"23\u00d788=2024\"},\n  @{old=\"64\u00d728=1792\"; new=\"52\u00d719=988\"},\n  @{old=\"73\u00d735=2555\"; new=\"50\u00d744=2200\"},\n  @{old=\"66\u00d747=3102\"; new=\"97\u00d759=5723\"},\n  @{old=\"78\u00d798=7644\"; new=\"97\u00d781=7857\"},\n  @{old=\"79\u00d765=5135\"; new=\"27\u00d794=2538\"},\n  @{old=\"20\u00d722=440\"; new=\"14\u00d755=770\"},\n  @{old=\"58\u00d724=1392\"; new=\"42\u00d746=1932\"},\n  @{old=\"92\u00d760=5520\"; new=\"72\u00d720=1440\"},\n  @{old=\"34\u00d731=1054\"; new=\"80\u00d734=2720\"},\n  @{old=\"65\u00d789=5785\"; new=\"82\u00d732=2624\"},\n  @{old=\"51\u00d773=3723\"; new=\"67\u00d755=3685\"},\n  @{old=\"36\u00d741=1476\"; new=\"26\u00d741=1066\"},\n  @{old=\"60\u00d737=2220\"; new=\"62\u00d798=6076\"},\n  @{old=\"13\u00d715=195\"; new=\"42\u00d775=3150\"},\n  @{old=\"78\u00d720=1560\"; new=\"17\u00d740=680\"},\n  @{old=\"29\u00d787=2523\"; new=\"50\u00d780=4000\"},\n)\n\nforeach ($p in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Execute($p.old, $false, $false, $false, $false, $false, $true, 1, $false, $p.new, 2)\n}"}
